$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$colA = " 🔎 : Joan"
$colB = " 📭 : Joan Martinez <joan_martinez.olivares@hotmail.com>"
$colC = "Joan Obtener Outlook para iOS<https://aka.ms/o0ukef> "
$colD = "📩 NUEVO 📩"

for ($r = 192; $r -le 194; $r++) {
    $ws.Cells.Item($r, 1).Value = $colA
    $ws.Cells.Item($r, 2).Value = $colB
    $ws.Cells.Item($r, 3).Value = $colC
    $ws.Cells.Item($r, 4).Value = $colD
}
